$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Text updates: "Ready for handoff" replaces "Handed back: in sync with en-US"
# and the handoff timestamps are bumped forward (report regenerated for
# handoff). These shared strings are referenced from all three sheets, so
# every referencing cell needs to be rewritten explicitly.
# ---------------------------------------------------------------------------

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-08-29 02:25:34"
$wsOverview.Range("G3").Value = "2016-08-29 02:25:34"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-08-29 02:25:30"
$wsZhCn.Range("H3").Value = "2016-08-29 02:25:30"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-08-29 02:25:34"
$wsDeDe.Range("H3").Value = "2016-08-29 02:25:34"

# ---------------------------------------------------------------------------
# Error Detail column (P) on the zh-cn and de-de sheets: the handback files
# are now flagged as stale versus the latest source.
# ---------------------------------------------------------------------------

$msg06af = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7a09e3a2e1bd0e193765430ca8f401529fbe0de8/e2e/06af091c-f622-4bc3-9bd1-2b7c9dacefbc.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/01a2589a919c0d66debe7de01e38179eb410ba5c/e2e/06af091c-f622-4bc3-9bd1-2b7c9dacefbc.md."
$msg36ec = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7a09e3a2e1bd0e193765430ca8f401529fbe0de8/e2e/36ec98e9-9065-443d-a93a-636bf5397cb3.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/01a2589a919c0d66debe7de01e38179eb410ba5c/e2e/36ec98e9-9065-443d-a93a-636bf5397cb3.md."

$wsZhCn.Range("P2").Value = $msg06af
$wsZhCn.Range("P3").Value = $msg36ec

$wsDeDe.Range("P2").Value = $msg06af
$wsDeDe.Range("P3").Value = $msg36ec

# ---------------------------------------------------------------------------
# Column width adjustments (narrower datetime columns, wider Error Detail).
# ColumnWidth is specified in characters; the underlying engine rounds to the
# nearest 1/6th of a character, so we pick the input nearest the recorded
# target so the stored width lands as close as possible.
# ---------------------------------------------------------------------------

$wsOverview.Columns.Item(5).ColumnWidth = 16.333333333333336   # E -> ~17.2159881591797
$wsOverview.Columns.Item(6).ColumnWidth = 16.333333333333336   # F -> ~17.2159881591797

$wsZhCn.Columns.Item(3).ColumnWidth = 16.333333333333336       # C -> ~17.2159881591797
$wsZhCn.Columns.Item(16).ColumnWidth = 39.16666666666667       # P -> 40

$wsDeDe.Columns.Item(3).ColumnWidth = 16.333333333333336       # C -> ~17.2159881591797
$wsDeDe.Columns.Item(16).ColumnWidth = 39.16666666666667       # P -> 40
